$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 338 (shifts rows 338-400 down to 339-401)
$ws.Rows.Item(338).Insert()

# Populate the newly inserted row 338 with the new record
$ws.Cells.Item(338, 1).Value = 4
$ws.Cells.Item(338, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(338, 3).Value = "Los Lagos"
$ws.Cells.Item(338, 4).Value = 44995
$ws.Cells.Item(338, 5).Value = 10
$ws.Cells.Item(338, 6).Value = 100112043
$ws.Cells.Item(338, 7).Value = "Pepino ensalada"
$ws.Cells.Item(338, 8).Value = "Sin especificar"
$ws.Cells.Item(338, 9).Value = "Primera"
$ws.Cells.Item(338, 10).Value = 400
$ws.Cells.Item(338, 11).Value = 14000
$ws.Cells.Item(338, 12).Value = 15000
$ws.Cells.Item(338, 13).Value = 14500
$ws.Cells.Item(338, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(338, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(338, 16).Value = 242
$ws.Cells.Item(338, 17).Value = 60
$ws.Cells.Item(338, 18).Value = "Hortaliza"
